$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" data column (J) to the table, mirroring the formatting
# already used for the preceding year column (I).
$ws.Range("I4:I14").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$ws.Range("J4").Value = 2021
$ws.Range("J5").Value = 24.4
$ws.Range("J6").Value = 45.7
$ws.Range("J7").Value = 38
$ws.Range("J8").Value = 51.3
$ws.Range("J9").Value = 51.5
$ws.Range("J10").Value = 13
$ws.Range("J11").Value = 36.4
$ws.Range("J12").Value = 27
$ws.Range("J13").Value = 2.7
$ws.Range("J14").Value = 40.4

# Row 3 (the blank separator row above the header) got a tighter height
# while editing.
$ws.Rows("3").RowHeight = 13.5

# Leave the selection where the editor finished up.
[void]$ws.Range("K18").Select()
